$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value2 = 1111.2572
$ws.Range("I15").Value2 = 1111.2572
$ws.Range("K15").Value2 = 3333.7716
$ws.Range("M15").Value2 = -3164.7716
# Row 132
$ws.Range("H132").Value2 = 1865.9678
$ws.Range("I132").Value2 = 1914.579
$ws.Range("K132").Value2 = 5743.737
$ws.Range("M132").Value2 = -3213.737
# Row 133
$ws.Range("H133").Value2 = 114555.2
$ws.Range("J133").Value2 = 114555.2
$ws.Range("L133").Value2 = 114555.2
$ws.Range("N133").Value2 = -124675.2
# Row 135
$ws.Range("H135").Value2 = 2800.375
$ws.Range("I135").Value2 = 557.5714
$ws.Range("K135").Value2 = 5018.1426
$ws.Range("M135").Value2 = -2483.1426
# Row 138
$ws.Range("H138").Value2 = 9555.777
$ws.Range("J138").Value2 = 9722.67
$ws.Range("L138").Value2 = 29168.01
$ws.Range("N138").Value2 = -39448.01

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value2 = 17846.018
$ws.Range("I32").Value2 = 17286.965
$ws.Range("J32").Value2 = 33499.5
$ws.Range("K32").Value2 = 17286.965
$ws.Range("L32").Value2 = 33499.5
$ws.Range("M32").Value2 = -16999.965
$ws.Range("N32").Value2 = -34073.5
# Row 61
$ws.Range("H61").Value2 = 7504015
$ws.Range("I61").Value2 = 10529674
$ws.Range("K61").Value2 = 10529674
$ws.Range("M61").Value2 = -10529462
# Row 74
$ws.Range("H74").Value2 = 2839.6316
$ws.Range("I74").Value2 = 2805.9443
$ws.Range("K74").Value2 = 2805.9443
$ws.Range("M74").Value2 = -1931.9443
# Row 77
$ws.Range("H77").Value2 = 2839.6316
$ws.Range("I77").Value2 = 2805.9443
$ws.Range("K77").Value2 = 14029.7215
$ws.Range("M77").Value2 = -9661.7215
# Row 97
$ws.Range("H97").Value2 = 982.6875
$ws.Range("I97").Value2 = 691.4516
$ws.Range("K97").Value2 = 691.4516
$ws.Range("M97").Value2 = -195.4516
# Row 110
$ws.Range("H110").Value2 = 9395.182000000001
$ws.Range("J110").Value2 = 6489.4
$ws.Range("L110").Value2 = 6489.4
$ws.Range("N110").Value2 = -10579.4
# Row 122
$ws.Range("H122").Value2 = 5993.3076
$ws.Range("I122").Value2 = 5873.04
$ws.Range("J122").Value2 = 9000
$ws.Range("K122").Value2 = 17619.12
$ws.Range("L122").Value2 = 27000
$ws.Range("M122").Value2 = -15169.12
$ws.Range("N122").Value2 = -31900
# Row 132
$ws.Range("H132").Value2 = 12503094
$ws.Range("I132").Value2 = 3536.5715
$ws.Range("K132").Value2 = 10609.7145
$ws.Range("M132").Value2 = -8079.7145
# Row 136
$ws.Range("H136").Value2 = 7504015
$ws.Range("I136").Value2 = 10529674
$ws.Range("K136").Value2 = 31589022
$ws.Range("M136").Value2 = -31586472
# Row 137
$ws.Range("H137").Value2 = 171168.1
$ws.Range("J137").Value2 = 171168.1
$ws.Range("L137").Value2 = 171168.1
$ws.Range("N137").Value2 = -181368.1

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value2 = 923297.8
$ws.Range("I86").Value2 = 2204447.8
$ws.Range("J86").Value2 = 8190.7144
$ws.Range("K86").Value2 = 2204447.8
$ws.Range("L86").Value2 = 8190.7144
$ws.Range("M86").Value2 = -2203324.8
$ws.Range("N86").Value2 = -10436.7144
# Row 89
$ws.Range("H89").Value2 = 923297.8
$ws.Range("I89").Value2 = 2204447.8
$ws.Range("J89").Value2 = 8190.7144
$ws.Range("K89").Value2 = 11022239
$ws.Range("L89").Value2 = 40953.572
$ws.Range("M89").Value2 = -11016623
$ws.Range("N89").Value2 = -52185.572
# Row 105
$ws.Range("H105").Value2 = 341185.8
$ws.Range("I105").Value2 = 430025.1
$ws.Range("K105").Value2 = 430025.1
$ws.Range("M105").Value2 = -428278.1
# Row 135
$ws.Range("H135").Value2 = 114995.25
$ws.Range("J135").Value2 = 114995.25
$ws.Range("L135").Value2 = 114995.25
$ws.Range("N135").Value2 = -125135.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value2 = 3677.3447
$ws.Range("I58").Value2 = 2850
$ws.Range("K58").Value2 = 2850
$ws.Range("M58").Value2 = -2647
# Row 92
$ws.Range("H92").Value2 = 9000
$ws.Range("J92").Value2 = 9000
$ws.Range("L92").Value2 = 9000
$ws.Range("N92").Value2 = -13992
# Row 103
$ws.Range("H103").Value2 = 43394.57
$ws.Range("I103").Value2 = 15396
$ws.Range("J103").Value2 = 54594
$ws.Range("K103").Value2 = 15396
$ws.Range("L103").Value2 = 54594
$ws.Range("M103").Value2 = -14224
$ws.Range("N103").Value2 = -56938
# Row 105
$ws.Range("H105").Value2 = 4994.5
$ws.Range("I105").Value2 = 1934.1765
$ws.Range("K105").Value2 = 1934.1765
$ws.Range("M105").Value2 = -187.1765
# Row 122
$ws.Range("H122").Value2 = 3549.75
$ws.Range("J122").Value2 = 1999
$ws.Range("L122").Value2 = 5997
$ws.Range("N122").Value2 = -10897
# Row 131
$ws.Range("H131").Value2 = 53999
$ws.Range("I131").Value2 = 0
$ws.Range("J131").Value2 = 53999
$ws.Range("K131").Value2 = 0
$ws.Range("L131").Value2 = ""
$ws.Range("M131").Value2 = 53999
$ws.Range("N131").Value2 = -64079
# Row 132
$ws.Range("H132").Value2 = 2952.0952
$ws.Range("I132").Value2 = 2741.7896
$ws.Range("K132").Value2 = 8225.3688
$ws.Range("M132").Value2 = -5695.3688
# Row 134
$ws.Range("H134").Value2 = 0
$ws.Range("I134").Value2 = 0
$ws.Range("J134").Value2 = 0
$ws.Range("K134").Value2 = 0
$ws.Range("L134").Value2 = ""
$ws.Range("M134").Value2 = ""
$ws.Range("N134").Value2 = 0
# Row 135
$ws.Range("H135").Value2 = 104998.07
$ws.Range("I135").Value2 = 40000
$ws.Range("J135").Value2 = 109997.92
$ws.Range("K135").Value2 = 40000
$ws.Range("L135").Value2 = 109997.92
$ws.Range("M135").Value2 = -34930
$ws.Range("N135").Value2 = -120137.92
# Row 136
$ws.Range("H136").Value2 = 3677.3447
$ws.Range("I136").Value2 = 2850
$ws.Range("K136").Value2 = 8550
$ws.Range("M136").Value2 = -6000
# Row 141
$ws.Range("H141").Value2 = 604827.7
$ws.Range("J141").Value2 = 604827.7
$ws.Range("L141").Value2 = 604827.7
$ws.Range("N141").Value2 = -615187.7

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value2 = 1012.0303
$ws.Range("I5").Value2 = 589.63635
$ws.Range("K5").Value2 = 1768.90905
$ws.Range("M5").Value2 = -1656.90905
# Row 12
$ws.Range("H12").Value2 = 40001936
$ws.Range("I12").Value2 = 100000936
$ws.Range("J12").Value2 = 2606.9333
$ws.Range("K12").Value2 = 300002808
$ws.Range("L12").Value2 = 7820.7999
$ws.Range("M12").Value2 = -300002635
$ws.Range("N12").Value2 = -8166.7999
# Row 76
$ws.Range("H76").Value2 = 21900
$ws.Range("J76").Value2 = 25000
$ws.Range("L76").Value2 = 75000
$ws.Range("N76").Value2 = -75766
# Row 79
$ws.Range("H79").Value2 = 21900
$ws.Range("J79").Value2 = 25000
$ws.Range("L79").Value2 = 75000
$ws.Range("N79").Value2 = -77652
# Row 135
$ws.Range("H135").Value2 = 1012.0303
$ws.Range("I135").Value2 = 589.63635
$ws.Range("K135").Value2 = 5306.72715
$ws.Range("M135").Value2 = -2771.72715
# Row 137
$ws.Range("H137").Value2 = 12418.305
$ws.Range("J137").Value2 = 15695.75
$ws.Range("L137").Value2 = 47087.25
$ws.Range("N137").Value2 = -57287.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 94
$ws.Range("H94").Value2 = 50000
$ws.Range("J94").Value2 = 50000
$ws.Range("L94").Value2 = 50000
$ws.Range("N94").Value2 = -51352
# Row 122
$ws.Range("H122").Value2 = 11014.866
$ws.Range("I122").Value2 = 7324.846
$ws.Range("J122").Value2 = 35000
$ws.Range("K122").Value2 = 21974.538
$ws.Range("L122").Value2 = 105000
$ws.Range("M122").Value2 = -19524.538
$ws.Range("N122").Value2 = -109900
# Row 136
$ws.Range("H136").Value2 = 58631.8
$ws.Range("J136").Value2 = 58631.8
$ws.Range("L136").Value2 = 175895.4
$ws.Range("N136").Value2 = -180995.4

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 14
$ws.Range("H14").Value2 = 9004.666999999999
$ws.Range("I14").Value2 = 9004
$ws.Range("J14").Value2 = 9005
$ws.Range("K14").Value2 = 9004
$ws.Range("L14").Value2 = 9005
$ws.Range("M14").Value2 = -8832
$ws.Range("N14").Value2 = -9349
# Row 21
$ws.Range("H21").Value2 = 6335.6665
$ws.Range("I21").Value2 = 4000
$ws.Range("K21").Value2 = 4000
$ws.Range("M21").Value2 = -3826
# Row 136
$ws.Range("H136").Value2 = 10365.417
$ws.Range("I136").Value2 = 13782.286
$ws.Range("K136").Value2 = 41346.858
$ws.Range("M136").Value2 = -38796.858

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 16
$ws.Range("H16").Value2 = 130899
$ws.Range("J16").Value2 = 130899
$ws.Range("L16").Value2 = 130899
$ws.Range("N16").Value2 = -131483
# Row 113
$ws.Range("H113").Value2 = 988.1667
$ws.Range("I113").Value2 = 843.3333
$ws.Range("J113").Value2 = 1133
$ws.Range("K113").Value2 = 2529.9999
$ws.Range("L113").Value2 = 3399
$ws.Range("M113").Value2 = -359.9998999999998
$ws.Range("N113").Value2 = -7739
# Row 132
$ws.Range("H132").Value2 = 3337262.8
$ws.Range("I132").Value2 = 5894
$ws.Range("K132").Value2 = 17682
$ws.Range("M132").Value2 = -15152
# Row 136
$ws.Range("H136").Value2 = 424086.88
$ws.Range("I136").Value2 = 7708.409
$ws.Range("J136").Value2 = 5004250
$ws.Range("K136").Value2 = 23125.227
$ws.Range("L136").Value2 = 15012750
$ws.Range("M136").Value2 = -20575.227
$ws.Range("N136").Value2 = -15017850

